$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# PayNowCC: dates refreshed, results stay "Pass"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("PayNowCC")
$ws.Range("B2").Value = "Wed Aug 14 22:43:32 EDT 2024"
$ws.Range("B3").Value = "Wed Aug 14 22:44:38 EDT 2024"
$ws.Range("B4").Value = "Wed Aug 14 22:45:34 EDT 2024"
$ws.Range("B5").Value = "Wed Aug 14 22:46:33 EDT 2024"

# ---------------------------------------------------------------------------
# PayNowCCSCF: row 2 date refreshed (stays Pass); rows 3-5 flip to Fail
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("PayNowCCSCF")
$ws.Range("B2").Value = "Wed Aug 14 22:47:30 EDT 2024"

$ws.Range("A3").Value = "Fail"
$ws.Range("B3").Value = "Wed Aug 14 22:48:36 EDT 2024"

$ws.Range("A4").Value = "Fail"
$ws.Range("B4").Value = "Wed Aug 14 22:53:09 EDT 2024"

$ws.Range("A5").Value = "Fail"
$ws.Range("B5").Value = "Wed Aug 14 22:54:32 EDT 2024"

# ---------------------------------------------------------------------------
# PayNowCCDCF: rows 2-5 flip to Fail with refreshed dates
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("PayNowCCDCF")

$ws.Range("A2").Value = "Fail"
$ws.Range("B2").Value = "Wed Aug 14 22:56:03 EDT 2024"

$ws.Range("A3").Value = "Fail"
$ws.Range("B3").Value = "Wed Aug 14 22:57:43 EDT 2024"

$ws.Range("A4").Value = "Fail"
$ws.Range("B4").Value = "Wed Aug 14 23:00:48 EDT 2024"

$ws.Range("A5").Value = "Fail"
$ws.Range("B5").Value = "Wed Aug 14 23:01:32 EDT 2024"

# ---------------------------------------------------------------------------
# OverUnderPay: rows 2-3 flip to Fail with refreshed dates
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("OverUnderPay")

$ws.Range("A2").Value = "Fail"
$ws.Range("B2").Value = "Wed Aug 14 23:05:47 EDT 2024"

$ws.Range("A3").Value = "Fail"
$ws.Range("B3").Value = "Wed Aug 14 23:29:09 EDT 2024"

# ---------------------------------------------------------------------------
# NoModifyAmount: row 2 flips to Fail with refreshed date
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("NoModifyAmount")
$ws.Range("A2").Value = "Fail"
$ws.Range("B2").Value = "Wed Aug 14 23:38:01 EDT 2024"

# ---------------------------------------------------------------------------
# NoOverPay: row 2 flips to Fail with refreshed date
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("NoOverPay")
$ws.Range("A2").Value = "Fail"
$ws.Range("B2").Value = "Thu Aug 15 03:22:29 EDT 2024"
